# Test env updates for verifying algorithm and input validation fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (project p1): shift 4th/5th choice values
$ws.Range("F2").Value = "s4"
$ws.Range("G2").Value = "s5"

# Row 3 (project p2): shift all choice values
$ws.Range("C3").Value = "s5"
$ws.Range("D3").Value = "s6"
$ws.Range("E3").Value = "s7"
$ws.Range("F3").Value = "s8"
$ws.Range("G3").Value = "s9"

# Row 4 (project p3): shift all choice values
$ws.Range("C4").Value = "s9"
$ws.Range("D4").Value = "s1"
$ws.Range("E4").Value = "s2"
$ws.Range("F4").Value = "s3"
$ws.Range("G4").Value = "s4"

# Update the active selection from G4 to G3
$ws.Range("G3").Select()
